$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.655.73"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.755.49"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "'324.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "'0.4582"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.75%  "

$ws.Range("D8").Value = "'0.3586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("D9").Value = "'0.07520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").Value = "'42.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").Value = "'1.097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").Value = "'20.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").Value = "'6.013"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.762.81"
$ws.Range("E15").Value = "  -1.95%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.103"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.97%  "

$ws.Range("D17").Value = "'92.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").Value = "'0.00001067"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("D19").Value = "'0.06414"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "'16.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("D22").Value = "'5.822"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("D23").Value = "27.705.61"
$ws.Range("E23").Value = "  -0.54%  "

$ws.Range("D24").Value = "'11.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("D25").Value = "'2.104"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("D26").Value = "'163.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.95%  "

$ws.Range("D27").Value = "'20.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("D28").Value = "1.959.40"
$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("D29").Value = "'2.080"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("D30").Value = "'126.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("E31").Value = "  -6.20%  "

$ws.Range("D32").Value = "'0.09193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.52%  "

$ws.Range("D33").Value = "'3.670"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").Value = "'5.540"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").Value = "'11.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("D36").Value = "'0.02298"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("D37").Value = "'0.2105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").Value = "'0.06058"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.982"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6337"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("D41").Value = "'1.207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.41%  "

$ws.Range("D42").Value = "'1.381"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("D43").Value = "'7.794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").Value = "'0.5913"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").Value = "'3.713"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").Value = "'123.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").Value = "'1.945"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("D49").Value = "'1.147"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.62%  "

$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("D51").Value = "'72.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.27%  "
